$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 281, shifting all subsequent rows down by one.
$ws.Rows("281:281").Insert()

# Populate the newly inserted row 281 with the new record.
$ws.Range("A281").Value = 4
$ws.Range("B281").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C281").Value = "Los Lagos"
$ws.Range("D281").Value = 44711
$ws.Range("E281").Value = 10
$ws.Range("F281").Value = 100112023
$ws.Range("G281").Value = "Brócoli"
$ws.Range("H281").Value = "Sin especificar"
$ws.Range("I281").Value = "Segunda"
$ws.Range("J281").Value = 750
$ws.Range("K281").Value = 1000
$ws.Range("L281").Value = 1000
$ws.Range("M281").Value = 1000
$ws.Range("N281").Value = '$/unidad'
$ws.Range("O281").Value = "Región del Maule"
$ws.Range("P281").Value = 1000
$ws.Range("Q281").Value = 1
$ws.Range("R281").Value = "Hortaliza"
